# FLT-Mileage Claim Form: stamp every NAV-bound content control with its
# alias (Title) / tag so the fields are identifiable by name in the UI,
# matching the data-binding xpath they are already wired to.

$tag = "#Nav: FLT_Mileage_Claim_Form/52122"

# --- main document body content controls (in document order) ---------
$d = $word.ActiveDocument

$bodyAliases = @(
    "#Nav: /Header/DepartmentCode",
    "#Nav: /Header/EmployeeName",
    "#Nav: /Header/Designation",
    "#Nav: /Header/Remarks",
    "#Nav: /Header/Lines/NumberofPassengers",
    "#Nav: /Header/Lines/Destination",
    "#Nav: /Header/Lines/VehicleModel",
    "#Nav: /Header/Lines/VehicleRegistrationNo_",
    "#Nav: /Header/Lines/EngineCapacity",
    "#Nav: /Header/Lines/TravelDate",
    "#Nav: /Header/TotalEstimatedMileage",
    "#Nav: /Header/ApprovedRatePerKm",
    "#Nav: /Header/Lines/ActualTotalCost",
    "#Nav: /Header/TransportOfficerSignature",
    "#Nav: /Header/TransportOfficersDate",
    "#Nav: /Header/ApproverSignature",
    "#Nav: /Header/ApproverDate"
)

for ($i = 1; $i -le $d.ContentControls.Count; $i++) {
    $cc = $d.ContentControls.Item($i)
    $cc.Title = $bodyAliases[$i - 1]
    $cc.Tag = $tag
}

# --- primary header content controls -----------------------------------
$headerAliases = @(
    "#Nav: /Header/CompanyPicture",
    "#Nav: /Header/CompanyName"
)

$hdr = $d.Sections.Item(1).Headers.Item(1)
for ($i = 1; $i -le $hdr.Range.ContentControls.Count; $i++) {
    $cc = $hdr.Range.ContentControls.Item($i)
    $cc.Title = $headerAliases[$i - 1]
    $cc.Tag = $tag
}

Write-Output "done"
